$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GRAY = 13421772  # 0x00CCCCCC (RGB cc,cc,cc)

# --- Legend rows (F column = label, centered; J column = swatch, left as-is) ---
$ws.Range("F2").Value = "Tillgänglig"
$ws.Range("F2").HorizontalAlignment = -4108  # xlCenter

$ws.Range("F3").Value = "Krock med icke obligatoriskt moment"
$ws.Range("F3").HorizontalAlignment = -4108

$ws.Range("F4").Value = "Krock med obligatoriskt moment"
$ws.Range("F4").HorizontalAlignment = -4108

# --- Week header rows (A column), gray fill + centered ---
$weeks = @(
    @{ Row = 7;  Label = "Vecka 3" },
    @{ Row = 13; Label = "Vecka 4" },
    @{ Row = 20; Label = "Vecka 5" },
    @{ Row = 27; Label = "Vecka 6" },
    @{ Row = 33; Label = "Vecka 7" },
    @{ Row = 39; Label = "Vecka 8" },
    @{ Row = 44; Label = "Vecka 9" },
    @{ Row = 48; Label = "Vecka 10" },
    @{ Row = 51; Label = "Vecka 13" },
    @{ Row = 54; Label = "Vecka 14" },
    @{ Row = 57; Label = "Vecka 15" },
    @{ Row = 60; Label = "Vecka 16" },
    @{ Row = 63; Label = "Vecka 17" },
    @{ Row = 66; Label = "Vecka 18" },
    @{ Row = 69; Label = "Vecka 19" }
)

foreach ($w in $weeks) {
    $cell = $ws.Cells.Item($w.Row, 1)
    $cell.Value = $w.Label
    $cell.HorizontalAlignment = -4108
    $cell.Interior.Color = $GRAY
}
